$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph that currently sits right after
#    the title (Heading1) paragraph.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2. Replace the final paragraph (the italic "Create a feature image..."
#    prompt) with two paragraphs:
#      - a new bold paragraph containing the page title text
#      - the same italic paragraph, but now containing the former meta
#        description text
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$finalPara = $d.Paragraphs($count)
$start = $finalPara.Range.Start
$end = $finalPara.Range.End   # include the paragraph mark so nothing is left behind
$targetRange = $d.Range($start, $end)

$newParagraphsXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Big Bot Crew Free: A Futuristic Slot Game | Review</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Join three tech-savvy teens in constructing a humanoid robot and win big with Big Bot Crew. Play for free and read our review for more details.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($newParagraphsXml)
